# Auto-generated Excel COM-interop script to update cached market-price
# values (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# Values mirror a scheduled refresh of live marketboard data; no formulas
# are involved anywhere in this workbook, so every cell is set directly.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 5666.3335
$ws.Range("I21").Value = 5666.3335
$ws.Range("K21").Value = 5666.3335
$ws.Range("M21").Value = -5198.3335
$ws.Range("H23").Value = 5666.3335
$ws.Range("I23").Value = 5666.3335
$ws.Range("K23").Value = 5666.3335
$ws.Range("M23").Value = -5432.3335
$ws.Range("H32").Value = 4285.769
$ws.Range("I32").Value = 3832
$ws.Range("J32").Value = 4569.375
$ws.Range("K32").Value = 3832
$ws.Range("L32").Value = 4569.375
$ws.Range("M32").Value = -3506
$ws.Range("N32").Value = -5221.375
$ws.Range("H33").Value = 167.57143
$ws.Range("I33").Value = 132.75
$ws.Range("K33").Value = 132.75
$ws.Range("M33").Value = 96.25
$ws.Range("H42").Value = 5362.143
$ws.Range("I42").Value = 2204.75
$ws.Range("J42").Value = 9572
$ws.Range("K42").Value = 6614.25
$ws.Range("L42").Value = 28716
$ws.Range("M42").Value = -6384.25
$ws.Range("N42").Value = -29176
$ws.Range("H43").Value = 2866.3333
$ws.Range("H62").Value = 4749.75
$ws.Range("I62").Value = 3666.3333
$ws.Range("K62").Value = 3666.3333
$ws.Range("M62").Value = -3042.3333
$ws.Range("H65").Value = 4749.75
$ws.Range("I65").Value = 3666.3333
$ws.Range("K65").Value = 18331.6665
$ws.Range("M65").Value = -15211.6665
$ws.Range("H74").Value = 112177.5
$ws.Range("I74").Value = 119845
$ws.Range("J74").Value = 12500
$ws.Range("K74").Value = 119845
$ws.Range("L74").Value = 12500
$ws.Range("M74").Value = -118909
$ws.Range("N74").Value = -14372
$ws.Range("H77").Value = 112177.5
$ws.Range("I77").Value = 119845
$ws.Range("J77").Value = 12500
$ws.Range("K77").Value = 599225
$ws.Range("L77").Value = 62500
$ws.Range("M77").Value = -594545
$ws.Range("N77").Value = -71860
$ws.Range("H112").Value = 3537.5
$ws.Range("J112").Value = 3516.6667
$ws.Range("L112").Value = 10550.0001
$ws.Range("N112").Value = -12766.0001
$ws.Range("H135").Value = 830.75
$ws.Range("I135").Value = 724.4545000000001
$ws.Range("K135").Value = 6520.0905
$ws.Range("M135").Value = -3985.0905
$ws.Range("H137").Value = 1545.5385
$ws.Range("I137").Value = 1562.909
$ws.Range("K137").Value = 4688.727000000001
$ws.Range("M137").Value = -2138.727000000001
$ws.Range("H141").Value = 3999.75
$ws.Range("J141").Value = 2999
$ws.Range("L141").Value = 8997
$ws.Range("N141").Value = -19357

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5358.1875
$ws.Range("I32").Value = 5358.1875
$ws.Range("K32").Value = 5358.1875
$ws.Range("M32").Value = -5071.1875
$ws.Range("H45").Value = 7866.5
$ws.Range("I45").Value = 8639.799999999999
$ws.Range("K45").Value = 8639.799999999999
$ws.Range("M45").Value = -8262.799999999999
$ws.Range("H56").Value = 39999
$ws.Range("J56").Value = 39999
$ws.Range("L56").Value = 39999
$ws.Range("N56").Value = -41483
$ws.Range("H61").Value = 4354.778
$ws.Range("I61").Value = 1299.75
$ws.Range("J61").Value = 6798.8
$ws.Range("K61").Value = 1299.75
$ws.Range("L61").Value = 6798.8
$ws.Range("M61").Value = -1087.75
$ws.Range("N61").Value = -7222.8
$ws.Range("H74").Value = 1458.8
$ws.Range("I74").Value = 1198.5
$ws.Range("K74").Value = 1198.5
$ws.Range("M74").Value = -324.5
$ws.Range("H77").Value = 1458.8
$ws.Range("I77").Value = 1198.5
$ws.Range("K77").Value = 5992.5
$ws.Range("M77").Value = -1624.5
$ws.Range("H132").Value = 2116
$ws.Range("I132").Value = 2116
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6348
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3818
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 4354.778
$ws.Range("I136").Value = 1299.75
$ws.Range("J136").Value = 6798.8
$ws.Range("K136").Value = 3899.25
$ws.Range("L136").Value = 20396.4
$ws.Range("M136").Value = -1349.25
$ws.Range("N136").Value = -25496.4

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1500
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1500
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 1500
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -1994
$ws.Range("H114").Value = 19998
$ws.Range("J114").Value = 19998
$ws.Range("L114").Value = 19998
$ws.Range("N114").Value = -28676
$ws.Range("H115").Value = 68390
$ws.Range("J115").Value = 68390
$ws.Range("L115").Value = 68390
$ws.Range("N115").Value = -71524
$ws.Range("H134").Value = 2499.5
$ws.Range("I134").Value = 2332.6667
$ws.Range("K134").Value = 6998.000100000001
$ws.Range("M134").Value = -4463.000100000001
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2162.75
$ws.Range("I31").Value = 1403.3334
$ws.Range("J31").Value = 4441
$ws.Range("K31").Value = 1403.3334
$ws.Range("L31").Value = 4441
$ws.Range("M31").Value = -1108.3334
$ws.Range("N31").Value = -5031
$ws.Range("H34").Value = 2162.75
$ws.Range("I34").Value = 1403.3334
$ws.Range("J34").Value = 4441
$ws.Range("K34").Value = 1403.3334
$ws.Range("L34").Value = 4441
$ws.Range("M34").Value = -1201.3334
$ws.Range("N34").Value = -4845
$ws.Range("H54").Value = 36449
$ws.Range("J54").Value = 33359.668
$ws.Range("L54").Value = 33359.668
$ws.Range("N54").Value = -34675.668
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H122").Value = 3475.2222
$ws.Range("I122").Value = 3475.2222
$ws.Range("K122").Value = 10425.6666
$ws.Range("M122").Value = -7975.6666
$ws.Range("H134").Value = 1805.25
$ws.Range("I134").Value = 1805.25
$ws.Range("K134").Value = 5415.75
$ws.Range("M134").Value = -2880.75
$ws.Range("H141").Value = 1000000
$ws.Range("J141").Value = 1000000
$ws.Range("L141").Value = 1000000
$ws.Range("N141").Value = -1010360

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10001851
$ws.Range("I4").Value = 22000034
$ws.Range("J4").Value = 3366
$ws.Range("K4").Value = 66000102
$ws.Range("L4").Value = 10098
$ws.Range("M4").Value = -65999990
$ws.Range("N4").Value = -10322
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H107").Value = 996
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 996
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2988
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -6828

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 3034
$ws.Range("J46").Value = 4046
$ws.Range("L46").Value = 4046
$ws.Range("N46").Value = -4358
$ws.Range("H113").Value = 1355.5
$ws.Range("J113").Value = 1355.5
$ws.Range("L113").Value = 1355.5
$ws.Range("N113").Value = -5695.5
$ws.Range("H126").Value = 3254.182
$ws.Range("J126").Value = 3325
$ws.Range("L126").Value = 9975
$ws.Range("N126").Value = -14915

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3120.5
$ws.Range("I40").Value = 2538.889
$ws.Range("J40").Value = 3868.2856
$ws.Range("K40").Value = 2538.889
$ws.Range("L40").Value = 3868.2856
$ws.Range("M40").Value = -2402.889
$ws.Range("N40").Value = -4140.2856
$ws.Range("H61").Value = 1714
$ws.Range("I61").Value = 877.8570999999999
$ws.Range("K61").Value = 877.8570999999999
$ws.Range("M61").Value = -675.8570999999999
$ws.Range("H113").Value = 1714
$ws.Range("I113").Value = 877.8570999999999
$ws.Range("K113").Value = 877.8570999999999
$ws.Range("M113").Value = 1292.1429
$ws.Range("H137").Value = 70000
$ws.Range("I137").Value = 30000
$ws.Range("J137").Value = 110000
$ws.Range("K137").Value = 30000
$ws.Range("L137").Value = 110000
$ws.Range("M137").Value = -24900
$ws.Range("N137").Value = -120200

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 44997.5
$ws.Range("J54").Value = 44997.5
$ws.Range("L54").Value = 44997.5
$ws.Range("N54").Value = -46037.5
$ws.Range("H96").Value = 1998
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H100").Value = 5264221
$ws.Range("I100").Value = 5883364.5
$ws.Range("K100").Value = 11766729
$ws.Range("M100").Value = -11766188
